$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) cells to match
# the latest scrape. The source values are stored as plain text (locale-
# formatted numbers such as "1.887.89" are not valid numeric literals, and
# the percentage cells keep their padding spaces), so each write forces the
# text interpretation via NumberFormat "@" and then restores the default
# "Normal" cell style so no stray number-format style is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.845.05'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.888.39'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7473'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -4.30%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.36'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.30'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07126'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08490'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.93%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7600'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.886.58'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.359'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.38'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.172'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.854.66'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.55'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007798'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.152.65'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.999'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.39%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1592'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.348'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.38'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.027'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.497'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.97%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.532'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.114'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.26%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05406'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.90%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.239'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7460'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.84%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01938'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.774'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4451'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.068'
$ws.Range('D42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.091.22'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.03%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.53'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.21%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8547'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.34'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.690'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.860'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.050'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.047.61'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.23%  '
$ws.Range('E51').Style = 'Normal'
